# CDS Study filter fixes
# Renames the "Cases" tab/query row to "Participants", fixes the duplicated
# lines in the StatQuery (column C) for every tab, bumps the query-box font
# size from 12 to 15, changes the selected cell, and resizes rows/columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix cell content
# ---------------------------------------------------------------------

# A2: "CasesTab" -> "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Column C (StatQuery) for all three data rows: remove the duplicated
# RETURN lines and replace with the corrected/trailing-space variant.
$fixedStatQuery = "MATCH (s:study)<--(p:participant)" + [char]10 + `
"OPTIONAL MATCH (p)<--(samp:sample)MATCH (s:study)<--(p:participant)" + [char]10 + `
"OPTIONAL MATCH (p)<--(samp:sample)" + [char]10 + `
"OPTIONAL MATCH (p)<--(diag:diagnosis)" + [char]10 + `
"OPTIONAL MATCH (samp)<--(f:file)" + [char]10 + `
"WITH DISTINCT samp,diag,s,p,f" + [char]10 + `
'WHERE s.study_name in ["University of Texas PDX Development and Trial Center Grant"]' + [char]10 + `
"RETURN" + [char]10 + `
"    count(distinct s) AS Studies," + [char]10 + `
"    count(distinct p) AS Participants," + [char]10 + `
"    count(distinct samp) AS Samples," + [char]10 + `
"    count(distinct f) AS ``Files``   "

$ws.Range("C2").Value = $fixedStatQuery
$ws.Range("C3").Value = $fixedStatQuery
$ws.Range("C4").Value = $fixedStatQuery

# ---------------------------------------------------------------------
# 2. Formatting: bump the query-box font to 15pt
# ---------------------------------------------------------------------
# NOTE: multi-area (comma separated) ranges only reliably apply to the
# first area in this host, so every contiguous block is set separately.

# Cells that keep their (no-wrap) look, just bigger font.
$ws.Range("A1:E1").Font.Size = 15
$ws.Range("A2").Font.Size = 15
$ws.Range("D2:E2").Font.Size = 15
$ws.Range("A3").Font.Size = 15
$ws.Range("D3:E3").Font.Size = 15
$ws.Range("A4").Font.Size = 15
$ws.Range("D4:E4").Font.Size = 15

# Cells that keep wrap-text on, just bigger font.
$ws.Range("B2:C4").Font.Size = 15
$ws.Range("B2:C4").WrapText = $true
$ws.Range("B5:C5").Font.Size = 15
$ws.Range("B5:C5").WrapText = $true
$ws.Range("C6").Font.Size = 15
$ws.Range("C6").WrapText = $true

# ---------------------------------------------------------------------
# 3. Row heights / column widths / view
# ---------------------------------------------------------------------

$ws.Rows.Item(4).RowHeight = 222.75

$ws.Columns.Item(1).ColumnWidth = 21.28515625
$ws.Columns.Item(2).ColumnWidth = 75.7109375
$ws.Columns.Item(3).ColumnWidth = 74.85546875
$ws.Columns.Item(4).ColumnWidth = 52
$ws.Columns.Item(5).ColumnWidth = 73.140625

$ws.Range("E9").Select()
